# Auto-generated edit script applying numeric corrections to the Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1883.4
$ws.Range("I4").Value = 1554.25
$ws.Range("J4").Value = 3200
$ws.Range("K4").Value = 1554.25
$ws.Range("L4").Value = 3200
$ws.Range("M4").Value = -1440.25
$ws.Range("N4").Value = -3428

$ws.Range("H20").Value = 1199
$ws.Range("I20").Value = 1199
$ws.Range("K20").Value = 1199
$ws.Range("M20").Value = -969

$ws.Range("H35").Value = 1199
$ws.Range("I35").Value = 1199
$ws.Range("K35").Value = 1199
$ws.Range("M35").Value = -820

$ws.Range("H38").Value = 1348.3334
$ws.Range("I38").Value = 1348.3334
$ws.Range("K38").Value = 4045.0002
$ws.Range("M38").Value = -3673.0002

$ws.Range("H39").Value = 282.4
$ws.Range("I39").Value = 282.4
$ws.Range("K39").Value = 847.1999999999999
$ws.Range("M39").Value = -551.1999999999999

$ws.Range("H40").Value = 2049.4375
$ws.Range("J40").Value = 1944.75
$ws.Range("L40").Value = 1944.75
$ws.Range("N40").Value = -2294.75

$ws.Range("H92").Value = 817
$ws.Range("I92").Value = 1036.75
$ws.Range("J92").Value = 597.25
$ws.Range("K92").Value = 1036.75
$ws.Range("L92").Value = 597.25
$ws.Range("M92").Value = 211.25
$ws.Range("N92").Value = -3093.25

$ws.Range("H112").Value = 3236.8462
$ws.Range("I112").Value = 696
$ws.Range("K112").Value = 2088
$ws.Range("M112").Value = -980

$ws.Range("H137").Value = 2615.6
$ws.Range("I137").Value = 2215.75
$ws.Range("K137").Value = 6647.25
$ws.Range("M137").Value = -4097.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4110
$ws.Range("I61").Value = 1679.625
$ws.Range("J61").Value = 7998.6
$ws.Range("K61").Value = 1679.625
$ws.Range("L61").Value = 7998.6
$ws.Range("M61").Value = -1467.625
$ws.Range("N61").Value = -8422.6

$ws.Range("H70").Value = 49999
$ws.Range("J70").Value = 49999
$ws.Range("L70").Value = 49999
$ws.Range("N70").Value = -50539

$ws.Range("H73").Value = 49999
$ws.Range("J73").Value = 49999
$ws.Range("L73").Value = 49999
$ws.Range("N73").Value = -51871

$ws.Range("H119").Value = 66332.664
$ws.Range("J119").Value = 66332.664
$ws.Range("L119").Value = 66332.664
$ws.Range("N119").Value = -76008.664

$ws.Range("H136").Value = 4110
$ws.Range("I136").Value = 1679.625
$ws.Range("J136").Value = 7998.6
$ws.Range("K136").Value = 5038.875
$ws.Range("L136").Value = 23995.8
$ws.Range("M136").Value = -2488.875
$ws.Range("N136").Value = -29095.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 706.86664
$ws.Range("I22").Value = 580.25
$ws.Range("J22").Value = 1213.3334
$ws.Range("K22").Value = 580.25
$ws.Range("L22").Value = 1213.3334
$ws.Range("M22").Value = -407.25
$ws.Range("N22").Value = -1559.3334

$ws.Range("H105").Value = 2283
$ws.Range("I105").Value = 2474.75
$ws.Range("J105").Value = 1899.5
$ws.Range("K105").Value = 2474.75
$ws.Range("L105").Value = 1899.5
$ws.Range("M105").Value = -727.75
$ws.Range("N105").Value = -5393.5

$ws.Range("H134").Value = 4805
$ws.Range("I134").Value = 4805
$ws.Range("K134").Value = 14415
$ws.Range("M134").Value = -11880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 672428.5600000001
$ws.Range("I6").Value = 783666.7
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 783666.7
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = -783553.7
$ws.Range("N6").Value = -5226

$ws.Range("H7").Value = 161.05556
$ws.Range("I7").Value = 126.6
$ws.Range("J7").Value = 333.33334
$ws.Range("K7").Value = 126.6
$ws.Range("L7").Value = 333.33334
$ws.Range("M7").Value = -13.59999999999999
$ws.Range("N7").Value = -559.33334

$ws.Range("H36").Value = 500
$ws.Range("I36").Value = 500
$ws.Range("K36").Value = 500
$ws.Range("M36").Value = -112

$ws.Range("H40").Value = 500
$ws.Range("I40").Value = 500
$ws.Range("K40").Value = 500
$ws.Range("M40").Value = -340

$ws.Range("H105").Value = 3311.5334
$ws.Range("I105").Value = 2587.8125
$ws.Range("K105").Value = 2587.8125
$ws.Range("M105").Value = -840.8125

$ws.Range("H106").Value = 17500
$ws.Range("J106").Value = 17500
$ws.Range("L106").Value = 17500
$ws.Range("N106").Value = -20024

$ws.Range("H107").Value = 767.4286
$ws.Range("I107").Value = 425.8889
$ws.Range("K107").Value = 425.8889
$ws.Range("M107").Value = 1494.1111

$ws.Range("H122").Value = 1830.25
$ws.Range("I122").Value = 1489.2858
$ws.Range("J122").Value = 2307.6
$ws.Range("K122").Value = 4467.857400000001
$ws.Range("L122").Value = 6922.799999999999
$ws.Range("M122").Value = -2017.857400000001
$ws.Range("N122").Value = -11822.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2478.2
$ws.Range("J51").Value = 2199.3333
$ws.Range("L51").Value = 6597.999899999999
$ws.Range("N51").Value = -7517.999899999999

$ws.Range("H58").Value = 1709.8
$ws.Range("I58").Value = 1029.5
$ws.Range("J58").Value = 2163.3333
$ws.Range("K58").Value = 3088.5
$ws.Range("L58").Value = 6489.999899999999
$ws.Range("M58").Value = -2960.5
$ws.Range("N58").Value = -6745.999899999999

$ws.Range("H94").Value = 11525.75
$ws.Range("I94").Value = 3301.5
$ws.Range("K94").Value = 9904.5
$ws.Range("M94").Value = -9228.5

$ws.Range("H119").Value = 3984.6667
$ws.Range("I119").Value = 3984.6667
$ws.Range("K119").Value = 11954.0001
$ws.Range("M119").Value = -7116.000100000001

$ws.Range("I137").Value = 2022.25
$ws.Range("J137").Value = 4634.3335
$ws.Range("K137").Value = 6066.75
$ws.Range("L137").Value = 13903.0005
$ws.Range("M137").Value = -966.75
$ws.Range("N137").Value = -24103.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 2933.3333
$ws.Range("I4").Value = 2933.3333
$ws.Range("K4").Value = 2933.3333
$ws.Range("M4").Value = -2821.3333

$ws.Range("H11").Value = 2620889.2
$ws.Range("I11").Value = 2944750.5
$ws.Range("K11").Value = 2944750.5
$ws.Range("M11").Value = -2944611.5

$ws.Range("H26").Value = 72499.5
$ws.Range("J26").Value = 72499.5
$ws.Range("L26").Value = 72499.5
$ws.Range("N26").Value = -73059.5

$ws.Range("H50").Value = 72499.5
$ws.Range("J50").Value = 72499.5
$ws.Range("L50").Value = 72499.5
$ws.Range("N50").Value = -73495.5

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("I70").Value = 7500
$ws.Range("K70").Value = 7500
$ws.Range("M70").Value = -7230

$ws.Range("I73").Value = 7500
$ws.Range("K73").Value = 7500
$ws.Range("M73").Value = -6564

$ws.Range("H80").Value = 3499
$ws.Range("I80").Value = 3499
$ws.Range("K80").Value = 3499
$ws.Range("M80").Value = -2501

$ws.Range("H83").Value = 3499
$ws.Range("I83").Value = 3499
$ws.Range("K83").Value = 17495
$ws.Range("M83").Value = -12503

$ws.Range("H113").Value = 1175.8889
$ws.Range("I113").Value = 1098
$ws.Range("J113").Value = 1448.5
$ws.Range("K113").Value = 1098
$ws.Range("L113").Value = 1448.5
$ws.Range("M113").Value = 1072
$ws.Range("N113").Value = -5788.5

$ws.Range("H132").Value = 2802.875
$ws.Range("I132").Value = 2802.875
$ws.Range("K132").Value = 8408.625
$ws.Range("M132").Value = -5878.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6934.8
$ws.Range("I7").Value = 3951.2222
$ws.Range("K7").Value = 3951.2222
$ws.Range("M7").Value = -3839.2222

$ws.Range("H16").Value = 407.42856
$ws.Range("I16").Value = 407.08334
$ws.Range("K16").Value = 407.08334
$ws.Range("M16").Value = -237.08334

$ws.Range("H22").Value = 942.2857
$ws.Range("I22").Value = 779.2
$ws.Range("K22").Value = 779.2
$ws.Range("M22").Value = -484.2

$ws.Range("H27").Value = 942.2857
$ws.Range("I27").Value = 779.2
$ws.Range("K27").Value = 779.2
$ws.Range("M27").Value = -672.2

$ws.Range("H40").Value = 4648.625
$ws.Range("I40").Value = 4000
$ws.Range("K40").Value = 4000
$ws.Range("M40").Value = -3864

$ws.Range("H122").Value = 6987.9
$ws.Range("I122").Value = 7755.7646
$ws.Range("J122").Value = 5983.769
$ws.Range("K122").Value = 23267.2938
$ws.Range("L122").Value = 17951.307
$ws.Range("M122").Value = -20817.2938
$ws.Range("N122").Value = -22851.307

$ws.Range("H126").Value = 6934.8
$ws.Range("I126").Value = 3951.2222
$ws.Range("K126").Value = 11853.6666
$ws.Range("M126").Value = -9383.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 15159
$ws.Range("I55").Value = 10596
$ws.Range("J55").Value = 16680
$ws.Range("K55").Value = 10596
$ws.Range("L55").Value = 16680
$ws.Range("M55").Value = -10319
$ws.Range("N55").Value = -17234

$ws.Range("H61").Value = 38495
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 38495
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 38495
$ws.Range("N61").Value = -39079
$ws.Range("M61").ClearContents()

$ws.Range("H123").Value = 106866.664
$ws.Range("J123").Value = 106866.664
$ws.Range("L123").Value = 106866.664
$ws.Range("N123").Value = -116666.664
